# Revert "Drop in files from RMI script"
# - Insert a new "Data Texas" worksheet between "About" and "HPPECbP" containing
#   commentary notes (and a hyperlink to the NREL study).
# - Change the hard-coded excess-capacity assumption on HPPECbP (cell B2, and
#   all the cells that derive from it) from 25% (0.25) to 10% (0.1).

$wb = $excel.ActiveWorkbook

$sheetAbout = $wb.Worksheets.Item("About")
$sheetHPP   = $wb.Worksheets.Item("HPPECbP")

$sheetAbout.Range("B19").Select() | Out-Null

# ------------------------------------------------------------------
# 1. Insert the new "Data Texas" sheet right after "About"
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheetAbout)
$newSheet.Name = "Data Texas"

$newSheet.Range("A1").Value = "Basically, they assume that electrolyzers will be sized to 125% of their average output. "
$newSheet.Range("A2").Value = "Or, put another way, they assume a capacity factor of 80%."
$newSheet.Range("B13").Value = "https://www.nrel.gov/docs/fy14osti/60528.pdf"
$newSheet.Hyperlinks.Add($newSheet.Range("B13"), "https://www.nrel.gov/docs/fy14osti/60528.pdf") | Out-Null
$newSheet.Range("A3").Value = "This seems like a weird spreadsheet. The study they cited is all about using electrolysis to produce hydrogen for fuel cell electric vehicles (FCEVs). The study's baseline scenario sizes"
$newSheet.Range("A4").Value = "the electrolyzers so that they can produce enough hydrogen to meet the FCEV demand. Then, as a sensitivity, they see how it would change things if the electrolyzers were oversized. "
$newSheet.Range("A5").Value = "They test them at 25% oversized and 50% oversized. In that case, the electrolyzers have excess capacity and can ramp up and down to provide flexible electricity consumption for the grid. "
$newSheet.Range("A7").Value = "The way EPS uses this data is confusing. For one, these results are meaningless for steam methane reforming and other hydrogen production processes that do not use electricity."
$newSheet.Range("A8").Value = "Otherwise, I assume they use these numbers to artificially inflate the capacity of hydrogen being installed which would drive up costs. "
$newSheet.Range("A10").Value = "In reality, hydrogen equipment will probably be used intensively. It is likely uneconomical to oversize electrolysis by 25% just to provide excess capacity for flexible grid services. Batteries are"
$newSheet.Range("A11").Value = "probably a much cheaper way to do this than electrolyzers, which are expensive. "
$newSheet.Range("A12").Value = "Moreover, an NREL study cited in another on of the EPS data sheets:"
$newSheet.Range("A14").Value = "assumes capacity factors of 90-98% for the different hydrogen producing technologies. So a capacity factor of 80%, which EPS assumes for this spreadsheet, seems quite conservative. "
$newSheet.Range("A16").Value = "I am updating their values from 25% to 10%, which is a capacity factor of 91% (still slightly conservative based on the NREL data)"
$newSheet.Range("A18").Value = "From a Texas standpoint, there is no reason to expect that Texas capacity factors would be lower or higher than US numbers. "

$newSheet.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 90
$newSheet.Range("C21").Select() | Out-Null

# ------------------------------------------------------------------
# 2. Update the excess-capacity assumption on HPPECbP from 25% to 10%
#    (all other cells in B2:AI6 are formulas referencing $B$2, so they
#    recalculate automatically). Re-fetch the worksheet reference since
#    the sheet collection changed after inserting "Data Texas".
# ------------------------------------------------------------------
$sheetHPP = $wb.Worksheets.Item("HPPECbP")
$sheetHPP.Range("B2").Value = 0.1
$sheetHPP.Activate() | Out-Null
$sheetHPP.Range("B3").Select() | Out-Null

Write-Output "done"
